# edit.ps1 - apply the "06-Features" commit:
#   1. Refresh the cached "datetimeFigureOut" date-field text that lives on
#      the slide master, every slide layout, the handout master and the
#      notes master (PowerPoint re-stamps these with the current date
#      whenever the deck is re-saved) from 10/11/22 -> 10/16/23.
#   2. Retitle slide 4 from "Statistical Features" to
#      "Statistical Features (STAT)".

$p = $ppt.ActivePresentation

$oldDate = "10/11/22"
$newDate = "10/16/23"

function Update-DateInShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a. Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DateInShapes $master.Shapes

# 1b. Every custom (slide) layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateInShapes $layout.Shapes
}

# 1c. Handout master - its placeholder shapes aren't directly editable via
# TextFrame.TextRange here, so go through the Headers/Footers dialog model.
$handout = $p.HandoutMaster
$handout.HeadersFooters.DateAndTime.Text = $newDate

# 1d. Notes master - same story as the handout master.
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = $newDate

# 2. Slide 4 title: "Statistical Features" -> "Statistical Features (STAT)".
$slide4 = $p.Slides.Item(4)
$title = $slide4.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Statistical Features (STAT)"
